$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '42.576.63'
$ws.Range('E2').Value = '  -1.18%  '
Set-TextValue $ws.Range('D3') '2.370.72'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').Value = '  -0.27%  '
Set-TextValue $ws.Range('D5') '333.89'
$ws.Range('E5').Value = '  +8.50%  '
Set-TextValue $ws.Range('D6') '101.11'
$ws.Range('E6').Value = '  -6.09%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.07%  '
Set-TextValue $ws.Range('D9') '0.638'
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('E10').Value = '  -4.84%  '
$ws.Range('E11').Value = '  -1.34%  '
$ws.Range('E12').Value = '  -3.90%  '
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('E14').Value = '  +0.32%  '
Set-TextValue $ws.Range('D15') '16.38'
$ws.Range('E15').Value = '  +0.21%  '
Set-TextValue $ws.Range('D16') '2.728.61'
$ws.Range('E16').Value = '  +0.48%  '
Set-TextValue $ws.Range('D17') '2.378.20'
$ws.Range('E17').Value = '  -3.19%  '
Set-TextValue $ws.Range('D18') '42.692.15'
$ws.Range('E18').Value = '  -0.89%  '
Set-TextValue $ws.Range('D19') '7.86'
$ws.Range('E19').Value = '  +7.67%  '
$ws.Range('E20').Value = '  -1.00%  '
Set-TextValue $ws.Range('D21') '3.82'
$ws.Range('E21').Value = '  +12.47%  '
Set-TextValue $ws.Range('D22') '75.79'
$ws.Range('E22').Value = '  +0.74%  '
Set-TextValue $ws.Range('D23') '271.22'
$ws.Range('E23').Value = '  +7.60%  '
Set-TextValue $ws.Range('D24') '2.33'
$ws.Range('E24').Value = '  -7.22%  '
Set-TextValue $ws.Range('D25') '10.06'
$ws.Range('E25').Value = '  +12.81%  '
$ws.Range('E26').Value = '  -0.03%  '
Set-TextValue $ws.Range('D27') '11.47'
$ws.Range('E27').Value = '  -3.82%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D28') '23.40'
$ws.Range('E28').Value = '  +3.98%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D29') '2.20'
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D30') '175.37'
$ws.Range('E30').Value = '  +1.33%  '
$ws.Range('B31').Value = 'WEMIXToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D31') '3.09'
$ws.Range('E31').Value = '  -2.57%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D32') '0.0911'
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D33') '35.56'
$ws.Range('E33').Value = '  -8.73%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D34') '6.13'
$ws.Range('E34').Value = '  +3.17%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D35') '0.134'
$ws.Range('E35').Value = '  +1.68%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D36') '4.62'
$ws.Range('E36').Value = '  -8.17%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D37') '0.0361'
$ws.Range('E37').Value = '  -4.26%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D38') '2.93'
$ws.Range('E38').Value = '  +8.16%  '
Set-TextValue $ws.Range('D39') '3.87'
$ws.Range('E39').Value = '  -4.83%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D40') '0.106'
$ws.Range('E40').Value = '  +3.47%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D41') '1.54'
$ws.Range('E41').Value = '  +3.34%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D42') '0.236'
$ws.Range('E42').Value = '  +2.53%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue $ws.Range('D43') '70.36'
$ws.Range('E43').Value = '  -2.10%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D44') '1.00'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D45') '118.44'
$ws.Range('E45').Value = '  +7.85%  '
$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue $ws.Range('D46') '89.68'
$ws.Range('E46').Value = '  +29.69%  '
Set-TextValue $ws.Range('D47') '12.10'
$ws.Range('E47').Value = '  -2.18%  '
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range('D48') '5.49'
$ws.Range('E48').Value = '  -1.85%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D49') '9.15'
$ws.Range('E49').Value = '  -1.91%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D50') '1.585.22'
$ws.Range('E50').Value = '  +5.90%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D51') '1.27'
$ws.Range('E51').Value = '  -1.25%  '
